$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - headers
$ws.Range("A1").Value = "Job_Id"
$ws.Range("B1").Value = "Job_Title"
$ws.Range("C1").Value = "Job_Description"
$ws.Range("D1").Value = "Total_Years_Min_Exp"
$ws.Range("E1").Value = "Total_Years_Max_Exp"
$ws.Range("F1").Value = "LinkedIn_Poster"
$ws.Range("G1").Value = "LinkedIn_Posted"
$ws.Range("H1").Value = "Resume_received"
$ws.Range("I1").Value = "Resume_downloaded"

# Row 2 - data
$ws.Range("A2").Value = "JD_001"
$ws.Range("B2").Value = "RPA Developer"
$ws.Range("C2").Value = "We are seeking a RPA Developer to design, develop, and support automation solutions.`nCollaborate with teams to streamline business processes using RPA tools like UiPath or Automation Anywhere. Join Akkodis to grow your skills in a dynamic, tech-driven environment"
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 4

$ws.Rows.Item(2).AutoFit()
